$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "sheet name:" $ws.Name
